$d = $word.ActiveDocument

# 1. Append " Is Baz" to the document title ("The Title" -> "The Title Is Baz"),
#    typed word-by-word / space-by-space so each word lands in its own run,
#    matching Word's normal typing behaviour.
$titleRange = $d.Paragraphs(1).Range
$titleRange.End = $titleRange.End - 1
$titleRange.InsertAfter(" ")
$titleRange.Collapse(0)
$titleRange.InsertAfter("Is")
$titleRange.Collapse(0)
$titleRange.InsertAfter(" ")
$titleRange.Collapse(0)
$titleRange.InsertAfter("Baz")

# 2. Update the running header that mirrors the title text.
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$hdr.Range.Find.Execute("Author / The Title / ", $false, $false, $false, $false, $false, $true, 1, $false, "Author / The Title Is Baz / ", 2)
